# Decisões Técnicas.docx - update table content per commit:
# "[Douglas Giordano] - Atualização decisões tecnicas."
#
# Summary of table changes (3-column table: N°, Questão Técnica, Decisão):
#   Row 05: "Banco de Dados" -> "SGBD"; "Mysql" text cleaned (bookmark removed
#           from middle of word, kept as single run wrapped in proofErr).
#   New row 06: "IDE BD" (+ tab stop) / "Mysql Workbench" inserted between the
#           old rows 05 and 06 (the old row 06 "SGBD"/"Mysql Workbench" is
#           replaced by this, since its text moved elsewhere).
#   Rows 07 "Ferramenta de Teste"/"JUnit" and 08 "Arquitetura do
#           Software"/"MVC + DAO" are unchanged (just shift down by one).
#   The two previously-empty trailing rows become:
#     Row 09: "Controle de Versionamento" / "GIT"
#     Row 10: "Metodologia de versionamento" / "GIT FLOW" (bookmark now sits
#             after the "10" run in the first cell).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-CellXml {
    param($cell, [string]$innerBodyXml)
    # InsertXML's behaviour depends on whether the target range is genuinely
    # empty: a cell with no real text (just the end-of-cell mark) has its
    # sole paragraph replaced in place (paragraph count unchanged); a cell
    # that already holds text gets the new paragraph(s) appended after the
    # old one(s), which then need deleting so only the fresh content stays.
    $oldCount = $cell.Range.Paragraphs.Count
    $full = $pkgHeader + '<w:body>' + $innerBodyXml + '</w:body>' + $pkgFooter
    $cell.Range.InsertXML($full)
    $newCount = $cell.Range.Paragraphs.Count
    $toDelete = $newCount - $oldCount
    for ($i = 1; $i -le $toDelete; $i++) {
        $cell.Range.Paragraphs.Item(1).Range.Delete()
    }
}

# ---- Row 05: "Banco de Dados" -> "SGBD" ----
$row05 = $t.Rows.Item(5)
Set-CellXml $row05.Cells.Item(2) '<w:p><w:r><w:t>SGBD</w:t></w:r></w:p>'

# ---- Row 05, 3rd cell: collapse "My" + bookmark + "sql" into one "Mysql" run ----
Set-CellXml $row05.Cells.Item(3) '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Mysql</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# ---- Insert new row 06 ("IDE BD") before the old row 06 ----
$oldRow06 = $t.Rows.Item(6)
$newRow06 = $t.Rows.Add($oldRow06)

Set-CellXml $newRow06.Cells.Item(1) '<w:p><w:r><w:t>06</w:t></w:r></w:p>'
Set-CellXml $newRow06.Cells.Item(2) '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1605"/></w:tabs></w:pPr><w:r><w:t>IDE BD</w:t></w:r><w:r><w:tab/></w:r></w:p>'
Set-CellXml $newRow06.Cells.Item(3) '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Mysql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Workbench</w:t></w:r></w:p>'

# The old row 06 ("SGBD"/"Mysql Workbench") is now redundant - its content has
# effectively moved into the new row above - so drop it. Rows "07" and "08"
# shift up into its place, unchanged.
$t.Rows.Item(7).Delete()

# ---- Fill the two previously-blank trailing rows ----
$row09 = $t.Rows.Item(9)
Set-CellXml $row09.Cells.Item(1) '<w:p><w:r><w:t>09</w:t></w:r></w:p>'
Set-CellXml $row09.Cells.Item(2) '<w:p><w:r><w:t>Controle de Versionamento</w:t></w:r></w:p>'
Set-CellXml $row09.Cells.Item(3) '<w:p><w:r><w:t>GIT</w:t></w:r></w:p>'

$row10 = $t.Rows.Item(10)
Set-CellXml $row10.Cells.Item(1) '<w:p><w:r><w:t>10</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Set-CellXml $row10.Cells.Item(2) '<w:p><w:r><w:t>Metodologia de versionamento</w:t></w:r></w:p>'
Set-CellXml $row10.Cells.Item(3) '<w:p><w:r><w:t>GIT FLOW</w:t></w:r></w:p>'

Write-Host "Done. Row count:" $t.Rows.Count
